$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hyperparameter configuration values (spatial clustering re-run).
# Each target cell is forced to Text format before the write so the numeric-
# looking strings are preserved as text (matching the source data), then the
# style is reset back to Normal so no stray formatting is introduced.
$updates = @{
    "C2" = "100"
    "D2" = "200"
    "E2" = "200"
    "F2" = "50"
    "C4" = "100"
    "D4" = "50"
    "C5" = "10"
    "D5" = "10"
    "E5" = "5"
    "F5" = "10"
    "C6" = "300"
    "F6" = "100"
    "C8" = "0.0001"
    "F8" = "1e-06"
    "C11" = "0.1"
    "C12" = "0.2"
    "D12" = "0.2"
    "E12" = "0.2"
    "C13" = "200"
    "D13" = "200"
    "E13" = "200"
    "D15" = "0.5"
    "E15" = "0.75"
    "C16" = "1000"
    "D16" = "200"
    "E16" = "100"
    "C17" = "10"
    "D17" = "100"
    "E17" = "100"
    "C18" = "4"
    "C19" = "0.01"
    "C20" = "28"
    "D20" = "18"
    "E20" = "28"
    "F20" = "8"
    "C21" = "100"
    "D21" = "100"
    "E21" = "200"
    "E22" = "4"
    "F22" = "5"
    "C23" = "0.05"
    "E23" = "0.05"
    "E24" = "10"
    "E25" = "0.1"
    "D26" = "[8]"
    "C27" = "600"
    "D27" = "600"
    "F28" = "10"
    "C29" = "4"
    "D29" = "4"
    "E29" = "5"
    "D30" = "6"
    "C31" = "140"
    "D31" = "120"
    "E31" = "100"
    "F31" = "120"
    "C32" = "0.001"
    "E32" = "0.001"
    "C33" = "1000"
    "E33" = "1000"
    "C34" = "0.1"
    "D34" = "0.1"
    "F34" = "0.1"
    "D35" = "800"
    "E35" = "600"
    "F35" = "600"
    "C36" = "5"
    "F36" = "4"
    "C37" = "0.02"
    "E37" = "0.01"
    "F37" = "0.01"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

Write-Host "Updated" $updates.Count "cells"